# "Added LM to Train Data"
# Insert a new "LM" column (Lm1..Lm20) between the "Sl.no" and "values"
# columns on the active sheet, and bold the new header cell to match the
# other (already-bold) header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "values"/"class" columns one place to the right,
# leaving a blank column B for the new data.
$ws.Columns("B:B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "LM"

# Fill Lm1..Lm20 for the 20 data rows (rows 2-21).
for ($i = 1; $i -le 20; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = "Lm$i"
}

# Match the bold styling already used on the other header cells.
$ws.Range("A1:B1").Font.Bold = $true
